$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("G2").Value = "2016-09-07 08:09:54"

$zhcn.Range("H2").Value = "2016-09-07 08:09:42"
$zhcn.Range("K2").Value = "2016-09-07 08:10:45"

$dede.Range("H2").Value = "2016-09-07 08:09:54"
$dede.Range("K2").Value = "2016-09-07 08:11:09"
